# bibliyor excel file loading
#
# Renames the "big  data" (double-space) sheet to "big data" (single
# space), widens two of its columns, and re-points the active
# window/selection at that sheet (it becomes the tab that is open/active
# instead of "statistics").

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename "big  data" -> "big data" ------------------------
# Renaming updates the <sheets> entry and the _xlnm._FilterDatabase
# defined name (which references 'big  data'!...) automatically.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "big data"

# --- Sheet 1: column width tweaks --------------------------------------
# Column A: 30.83203125 -> 32.83203125 (one unit wider)
# Column S (19th): 7.1640625 -> 12.1640625 (five units wider)
# ColumnWidth is expressed in "characters"; the stored OOXML <col width>
# is characters + 5/7, so subtract 5/7 before assigning.
$ws1.Columns.Item(1).ColumnWidth = 32.83203125 - (5/7)
$ws1.Columns.Item(19).ColumnWidth = 12.1640625 - (5/7)

# --- Window/selection: sheet 1 becomes the active tab ------------------
# Previously "statistics" (sheet 3) was the active/selected tab with
# cell D5 selected on "big  data" and B8 selected on "statistics".
# Now "big data" (sheet 1) is the active tab, scrolled down so row 35 is
# at the top, with C2 selected; "statistics" is no longer the active tab.
$ws1.Activate()

$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1

$ws1.Range("C2").Select()
